## Timelog.xlsx update
## - Adjust the Time Worked formula (E column) to handle shifts that cross
##   midnight (End < Start) by wrapping to the next day.
## - Fill in the logged entries for week 1 (rows 5-13): date, start/end time
##   and a short summary of the work done.
## - Switch the "Total Time Worked" cell (H4) to an elapsed-time format since
##   the total can now exceed 24 hours.
## - Move the active selection/scroll position to reflect where the user was
##   last working (around F14).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

## --- Fix the Time Worked formula so it wraps past midnight -----------------
## E4 is a standalone formula; E5:E68 and E69:E100 are each one shared-formula
## block, so re-entering the formula across each whole block keeps them as a
## single shared formula (matching how Excel itself would fill it down).
$ws.Range("E4").Formula = "=IF(D4-C4<0, D4-C4+1, D4-C4)"
$ws.Range("E5:E68").Formula = "=IF(D5-C5<0, D5-C5+1, D5-C5)"
$ws.Range("E69:E100").Formula = "=IF(D69-C69<0, D69-C69+1, D69-C69)"

## --- Fill in the week's time log entries ------------------------------------
$entries = @(
    @{ Row = 5;  Date = 45195; Start = 0.625;                Finish = 0.65763888888888888; Summary = "Weekly meeting" },
    @{ Row = 6;  Date = 45196; Start = 0.16666666666666666;  Finish = 0.21041666666666667;  Summary = "Updated meeting minutes, did some more research" },
    @{ Row = 7;  Date = 45196; Start = 0.9277777777777777;   Finish = 0.12986111111111112;  Summary = "Researched different graphics and plotting libraries, and example use cases" },
    @{ Row = 8;  Date = 45200; Start = 0.041666666666666664; Finish = 0.25625000000000003;  Summary = "Experimented with headless OpenGL, researched libraries, prepared more test data" },
    @{ Row = 9;  Date = 45202; Start = 0.17013888888888887;  Finish = 0.28125;               Summary = "Prepared slides, continued research into Widgets and H264 streaming" },
    @{ Row = 10; Date = 45202; Start = 0.5416666666666666;   Finish = 0.5625;                Summary = "Weekly meeting" },
    @{ Row = 11; Date = 45204; Start = 0.1125;                Finish = 0.1875;                Summary = "Meeting minutes" },
    @{ Row = 12; Date = 45206; Start = 0.3340277777777778;   Finish = 0.4548611111111111;    Summary = "Integrated Jupyter widget template" },
    @{ Row = 13; Date = 45209; Start = 0.24930555555555556;  Finish = 0.48541666666666666;   Summary = "Implemented prototype of the Jupyter widget which integrates an OpenGL renderer" }
)

foreach ($e in $entries) {
    $r = $e.Row
    $ws.Cells.Item($r, 2).Value = $e.Date      # B: Date
    $ws.Cells.Item($r, 3).Value = $e.Start     # C: Time Start
    $ws.Cells.Item($r, 4).Value = $e.Finish    # D: Time End
    $ws.Cells.Item($r, 6).Value = $e.Summary   # F: Summary
}

## --- Total Time Worked can now exceed 24 hours: use an elapsed-time format -
$ws.Range("H4").NumberFormat = "[h]:mm:ss"

## --- Tidy up the bottom border under the last data row ----------------------
$ws.Range("E100").Borders.Item(9).LineStyle = -4142   # xlEdgeBottom -> xlLineStyleNone

## --- Restore the view to where the user left off ---------------------------
$ws.Range("F14").Select()
$excel.ActiveWindow.ScrollRow = 3
$excel.ActiveWindow.ScrollColumn = 1
